$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet (workbook.xml <sheet name="..."/>)
$ws.Name = "thinBasic_Test_58602_255"

# New set of random values for column B (rows 1-20); column D holds
# =SIN(Bn) shared formulas that recalc automatically, and B22/D22 hold
# SUM() totals that recalc automatically as well.
$values = @(1618, 1594, 467, 757, 1520, 154, 1047, 1745, 1873, 1831, 1785, 809, 824, 1300, 1205, 384, 1421, 1471, 1511, 1006)

# Column B is formatted as Text (numFmt "@"), so a plain .Value assignment
# would be stored as text (matching real Excel behavior for text-formatted
# cells). Temporarily switch to General so the values land as numbers, the
# way they were originally, then restore the original number format.
$rng = $ws.Range("B1:B20")
$originalFormat = $rng.NumberFormat
$rng.NumberFormat = "General"

for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item($i + 1, 2).Value = $values[$i]
}

$rng.NumberFormat = $originalFormat
